$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.303.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.13%  "
$ws.Range("D3").Value = "'1.868.28"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'234.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").Value = "'0.4702"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  -1.39%  "
$ws.Range("D9").Value = "'0.06577"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("D11").Value = "'0.07816"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "'96.86"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.85%  "
$ws.Range("D13").Value = "'1.871.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "'0.6964"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "'5.092"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'268.53"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "'30.265.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "'13.78"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").Value = "'0.000007684"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("E20").Value = "  -0.03%  "
$ws.Range("D21").Value = "'2.110.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'5.247"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "'6.168"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").Value = "'9.539"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("D26").Value = "'166.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.52%  "
$ws.Range("D27").Value = "'18.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'1.941"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "'1.363"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "
$ws.Range("D30").Value = "'0.09921"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.40%  "
$ws.Range("D31").Value = "'4.365"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.32%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("D33").Value = "'4.053"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").Value = "'0.04736"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("D35").Value = "'1.132"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.22%  "
$ws.Range("D36").Value = "'0.7042"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("E37").Value = "  +0.35%  "
$ws.Range("D38").Value = "'0.01877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("D39").Value = "'2.770"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").Value = "'6.335"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("D41").Value = "'72.71"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.66%  "
$ws.Range("D42").Value = "'1.953"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.67%  "
$ws.Range("D43").Value = "'0.4177"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "'0.8362"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.73%  "
$ws.Range("D46").Value = "'103.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "'967.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.07%  "
$ws.Range("D48").Value = "'7.116"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.59%  "
$ws.Range("D49").Value = "'9.124"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").Value = "'34.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.04%  "
$ws.Range("D51").Value = "'0.05682"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.31%  "
